# Update "paises.xlsx" - countries & provincias Spain data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "data as of" timestamp in A1 (09:35 -> 10:05)
$ws.Range("A1").Value = "Datos actualizados a 25 de Mayo de 2020 a las 10:05"

# Row 6: Rusia - refreshed case counts (no reordering)
$ws.Range("B6").Value = 353427
$ws.Range("C6").Value = 8946
$ws.Range("D6").Value = 118798
$ws.Range("E6").Value = 230996
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 92
$ws.Range("H6").Value = 3633

# Rows 50-52: refreshed data causes a re-sort so that Afganistan moves above
# Serbia/Panama. Row 50 becomes Afganistan (new data), row 51 becomes the
# country previously on row 50 (Serbia), row 52 becomes the country
# previously on row 51 (Panama).
$ws.Range("A50").Value = "Afganistan"
$ws.Range("B50").Value = 11173
$ws.Range("C50").Value = 591
$ws.Range("D50").Value = 1097
$ws.Range("E50").Value = 9857
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 1
$ws.Range("H50").Value = 219

$ws.Range("A51").Value = "Serbia"
$ws.Range("B51").Value = 11159
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 5857
$ws.Range("E51").Value = 5064
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 238

$ws.Range("A52").Value = "Panama"
$ws.Range("B52").Value = 10926
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 6279
$ws.Range("E52").Value = 4341
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 306

# Row 94: Cuba - refreshed case counts (no reordering)
$ws.Range("B94").Value = 1635
$ws.Range("C94").Value = 12
$ws.Range("D94").Value = 1138
$ws.Range("E94").Value = 434
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 63

# Row 97: Islandia - refreshed case counts (no reordering)
$ws.Range("B97").Value = 1511
$ws.Range("C97").Value = 2
$ws.Range("D97").Value = 1307
$ws.Range("E97").Value = 176
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 28

# Rows 154-155: refreshed data causes Uganda to move above Birmania.
# Row 154 becomes Uganda (new data), row 155 becomes the country
# previously on row 154 (Birmania).
$ws.Range("A154").Value = "Uganda"
$ws.Range("B154").Value = 212
$ws.Range("C154").Value = 14
$ws.Range("D154").Value = 68
$ws.Range("E154").Value = 144
$ws.Range("F154").Value = 0
$ws.Range("G154").Value = 0
$ws.Range("H154").Value = 0

$ws.Range("A155").Value = "Birmania"
$ws.Range("B155").Value = 201
$ws.Range("C155").Value = 0
$ws.Range("D155").Value = 122
$ws.Range("E155").Value = 73
$ws.Range("F155").Value = 0
$ws.Range("G155").Value = 0
$ws.Range("H155").Value = 6
